$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells remain text, matching the workbook's inline-string format
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.167.39"
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").Value = "1.853.42"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "237.47"
$ws.Range("E5").Value = "  -2.07%  "

$ws.Range("E6").Value = "  -4.95%  "

$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3041"
$ws.Range("E8").Value = "  -3.48%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07610"
$ws.Range("E9").Value = "  +3.22%  "

$ws.Range("D10").Value = "23.11"
$ws.Range("E10").Value = "  -6.05%  "

$ws.Range("D11").Value = "0.08115"
$ws.Range("E11").Value = "  -1.04%  "

$ws.Range("D12").Value = "1.969.46"
$ws.Range("E12").Value = "  +4.18%  "

$ws.Range("D13").Value = "0.7222"
$ws.Range("E13").Value = "  -3.13%  "

$ws.Range("D14").Value = "5.182"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("D15").Value = "89.40"
$ws.Range("E15").Value = "  -3.42%  "

$ws.Range("D16").Value = "29.152.87"
$ws.Range("E16").Value = "  -2.14%  "

$ws.Range("D17").Value = "0.000007791"
$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("D18").Value = "5.713"
$ws.Range("E18").Value = "  -4.98%  "

$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  -2.04%  "

$ws.Range("D20").Value = "233.67"
$ws.Range("E20").Value = "  -5.28%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "2.098.76"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "7.406"
$ws.Range("E24").Value = "  -4.14%  "

$ws.Range("D25").Value = "161.49"
$ws.Range("E25").Value = "  -1.58%  "

$ws.Range("D26").Value = "8.930"
$ws.Range("E26").Value = "  -3.38%  "

$ws.Range("D27").Value = "0.1427"
$ws.Range("E27").Value = "  -5.22%  "

$ws.Range("D28").Value = "18.00"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").Value = "1.954"
$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").Value = "4.501"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").Value = "1.482"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("D33").Value = "3.999"
$ws.Range("E33").Value = "  -4.09%  "

$ws.Range("D34").Value = "0.05153"
$ws.Range("E34").Value = "  -5.87%  "

$ws.Range("D35").Value = "1.182"
$ws.Range("E35").Value = "  -4.00%  "

$ws.Range("D36").Value = "0.7044"
$ws.Range("E36").Value = "  -4.04%  "

$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("D38").Value = "2.674"
$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").Value = "0.01848"
$ws.Range("E39").Value = "  -3.38%  "

$ws.Range("D40").Value = "2.680"
$ws.Range("E40").Value = "  -2.21%  "

$ws.Range("D41").Value = "0.9108"
$ws.Range("E41").Value = "  +2.02%  "

$ws.Range("D42").Value = "1.103.79"
$ws.Range("E42").Value = "  +6.23%  "

$ws.Range("E43").Value = "  -0.80%  "

$ws.Range("D44").Value = "0.4271"
$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("D45").Value = "70.38"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("D47").Value = "101.82"
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("D48").Value = "1.775"
$ws.Range("E48").Value = "  -1.94%  "

$ws.Range("D49").Value = "1.993.71"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").Value = "9.137"
$ws.Range("E50").Value = "  -4.90%  "

$ws.Range("D51").Value = "6.980"
$ws.Range("E51").Value = "  -6.57%  "
